$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Formatting baseline: replicate row 6's look-and-feel down through the
#    six new rows (7-12) before filling in content, matching the existing
#    table's visual pattern.
# ---------------------------------------------------------------------------
$ws.Range("A6:H6").Copy()
$ws.Range("A7:H12").PasteSpecial(-4122)

# F8 / F10 are "Read" (teal) entries like F2, not "Imp" (purple) like F6.
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# F11 / F12 are "V Imp" entries rendered with a solid red fill; reset to the
# plain bordered look first, then paint it red.
$ws.Range("A2").Copy()
$ws.Range("F11:F12").PasteSpecial(-4122)
$ws.Range("F11:F12").Interior.Color = 255
$ws.Range("F11:F12").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Row content. Filled in row-major, left-to-right order so that new
#    shared-string entries are appended in the same order the author's
#    Excel session would have created them.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Move all Zeros to the end of the array"
$ws.Range("E7").Value = "CN"
$ws.Range("F7").Value = "Imp"
$ws.Range("G7").Value = "Pass"
$ws.Range("H7").Value = "22-12-2023"

$ws.Range("B8").Value = "Union of Two Sorted Arrays"
$ws.Range("E8").Value = "CN"
$ws.Range("F8").Value = "Read"
$ws.Range("G8").Value = "Pass"
$ws.Range("H8").Value = "22-12-2023"

$ws.Range("B9").Value = "Count Maximum Consecutive One's in the array"
$ws.Range("E9").Value = "LC"
$ws.Range("F9").Value = "Imp"
$ws.Range("G9").Value = "Pass"
$ws.Range("H9").Value = "22-12-2023"

$ws.Range("B10").Value = "Find the number that appears once, and the other numbers twice"
$ws.Range("E10").Value = "CN"
$ws.Range("F10").Value = "Read"
$ws.Range("G10").Value = "Pass"
$ws.Range("H10").Value = "22-12-2023"

$ws.Range("B11").Value = "Longest Subarray with given Sum K(Positives)"
$ws.Range("E11").Value = "CN"
$ws.Range("F11").Value = "V Imp"
$ws.Range("G11").Value = "Pass"
$ws.Range("H11").Value = "24-12-2023"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Longest Subarray with sum K | [Postives and Negatives]"
$ws.Range("E12").Value = "CN"
$ws.Range("F12").Value = "V Imp"
$ws.Range("G12").Value = "Pass"
$ws.Range("H12").Value = "24-12-2023"

# ---------------------------------------------------------------------------
# 3. Hyperlinks for the newly added problems.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.codingninjas.com/studio/problems/move-all-zeros-to-the-end-of-the-array_1083479?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.codingninjas.com/studio/problems/union-of-two-sorted-arrays_1266479?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://leetcode.com/problems/max-consecutive-ones/")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.codingninjas.com/studio/problems/find-missing-repeating-numbers_6828164?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.codingninjas.com/studio/problems/longest-subarray-with-sum-k-_6763156?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.codingninjas.com/studio/problems/longest-subarray-with-sum-k_7426776?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")

# Adding a hyperlink re-paints its cell with Excel's built-in "Hyperlink"
# look; re-apply the sheet's own hyperlink-cell format (border + style,
# matching B2:B6) so B7:B12 stay visually consistent with the rest of the
# table instead of picking up a second, slightly-different hyperlink xf.
$ws.Range("B2").Copy()
$ws.Range("B7:B12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Column width + selection to mirror the saved UI state.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 60

$ws.Range("C10").Select()
